$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.189.82'
$ws.Range("E2").Value = '  -0.32%  '
$ws.Range("D3").Value = '1.913.76'
$ws.Range("E3").Value = '  +0.13%  '
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = "'325.61"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.66%  '
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.21%  '
$ws.Range("D7").Value = "'0.4616"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = "'0.3895"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.42%  '
$ws.Range("D9").Value = "'0.07870"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.67%  '
$ws.Range("D10").Value = "'0.9925"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.88%  '
$ws.Range("D11").Value = "'22.01"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.67%  '
$ws.Range("D12").Value = '1.924.75'
$ws.Range("E12").Value = '  +0.73%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = "'5.763"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.03%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = "'7.050"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.57%  '
$ws.Range("D15").Value = "'0.07049"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.13%  '
$ws.Range("D16").Value = "'88.08"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.69%  '
$ws.Range("D17").Value = "'1.003"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.10%  '
$ws.Range("D18").Value = "'0.000009962"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.93%  '
$ws.Range("D19").Value = "'17.09"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.23%  '
$ws.Range("D20").Value = "'1.001"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.22%  '
$ws.Range("D21").Value = '29.193.73'
$ws.Range("E21").Value = '  -0.28%  '
$ws.Range("D22").Value = "'5.337"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("E23").Value = '  +0.37%  '
$ws.Range("D24").Value = '2.130.73'
$ws.Range("E24").Value = '  -4.67%  '
$ws.Range("D25").Value = "'2.091"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.11%  '
$ws.Range("D26").Value = "'156.13"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.47%  '
$ws.Range("D27").Value = "'19.47"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("D28").Value = "'5.893"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.78%  '
$ws.Range("D29").Value = "'118.88"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.41%  '
$ws.Range("D30").Value = "'1.870"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -5.39%  '
$ws.Range("D31").Value = "'0.09338"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.55%  '
$ws.Range("D32").Value = "'0.8929"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -3.30%  '
$ws.Range("D33").Value = "'5.223"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.39%  '
$ws.Range("D34").Value = "'1.321"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -2.69%  '
$ws.Range("D35").Value = "'3.134"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -4.32%  '
$ws.Range("D36").Value = "'0.05783"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.98%  '
$ws.Range("E37").Value = '  -2.47%  '
$ws.Range("D38").Value = "'0.02090"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.79%  '
$ws.Range("D39").Value = "'0.9999"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.27%  '
$ws.Range("D40").Value = "'0.5704"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -0.82%  '
$ws.Range("D41").Value = "'7.668"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -3.23%  '
$ws.Range("D42").Value = "'0.1806"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.32%  '
$ws.Range("D43").Value = "'9.744"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.14%  '
$ws.Range("D44").Value = "'0.000002853"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +77.03%  '
$ws.Range("D45").Value = "'11.91"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.69%  '
$ws.Range("D46").Value = "'0.5355"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.17%  '
$ws.Range("D47").Value = "'2.197"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -4.54%  '
$ws.Range("D48").Value = "'0.06968"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.31%  '
$ws.Range("D49").Value = "'1.843"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.79%  '
$ws.Range("D50").Value = "'2.546"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.31%  '
$ws.Range("D51").Value = "'112.78"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.59%  '
